$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.218.58"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.687.41"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0892"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.926.20"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "1.697.98"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.556"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "27.214.58"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.40%  "
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").Value = "1.575.01"
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  +5.73%  "
$ws.Range("E37").Value = "  +3.47%  "
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0175"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("E41").Value = "  +3.49%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("D45").Value = "1.833.98"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("E48").Value = "  +5.55%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.76%  "
